# The presentation shipped with two theme parts:
#   ppt/theme/theme1.xml  ("Integral"     colour scheme) -> used by the slide master
#   ppt/theme/theme2.xml  ("Office Theme" colour scheme) -> used by the notes master
#
# The commit swaps the contents of those two theme parts (theme1.xml becomes
# the former "Office Theme" content, theme2.xml becomes the former "Integral"
# content) - the fontScheme/fmtScheme blocks are byte-identical between the
# two themes, only the 12 clrScheme colours (and the cosmetic name=
# attributes) actually differ.
#
# This host's object model exposes a single live Theme/ThemeColorScheme
# (reachable from the SlideMaster, the NotesMaster, individual Slides, ...)
# that is backed by ppt/theme/theme1.xml, so we reproduce the colour part of
# the swap by writing the "Office Theme" palette into that shared theme
# through the standard ThemeColorScheme COM surface.

function ToOle($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# Target palette = the "Office Theme" colours that used to live in theme2.xml,
# in MsoThemeColorSchemeIndex order (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = ToOle($officeThemeColors[$i - 1])
}
